$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cross-period average comparison formulas
$ws.Range("S4").Formula = "=AVERAGE(O2:O5)+AVERAGE(O15:O25)"
$ws.Range("S5").Formula = "=AVERAGE(O6:O14)+AVERAGE(O26:O34)"
$ws.Range("T5").Formula = "=S5-S4"

# Move summary formulas from column Q to column P, converting the first two to live formulas
$ws.Range("P36").Formula = "=AVERAGE(O2:O5,O15:O25)"
$ws.Range("Q36").ClearContents()

$ws.Range("P37").Formula = "=AVERAGE(O26:O34,O6:O14)"
$ws.Range("Q37").ClearContents()

$ws.Range("P38").Value = 0.97499999999999998
$ws.Range("Q38").ClearContents()

# Update view: scroll + selection
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("P38").Select()
